# Weekly update: a new price-report row for Puerro (Vega Central Mapocho de
# Santiago) is inserted at row 19, pushing all the existing weekly records
# (previously rows 19-80) down by one row. The oldest record (previously row
# 80) ends up as the new last row (81).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19; everything from row 19 down (through the
# former last row, 80) shifts down to make room, landing on 20-81.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with this week's data point.
$ws.Range("A19").Value = 9
$ws.Range("B19").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C19").Value = "Metropolitana"
$ws.Range("D19").Value = 44526
$ws.Range("E19").Value = 13
$ws.Range("F19").Value = 100112005
$ws.Range("G19").Value = "Puerro"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 170
$ws.Range("K19").Value = 6000
$ws.Range("L19").Value = 7000
$ws.Range("M19").Value = 6500
$ws.Range("N19").Value = "$/paquete 20 unidades"
$ws.Range("O19").Value = "Provincia de Chacabuco"
$ws.Range("P19").Value = 325
$ws.Range("Q19").Value = 20
$ws.Range("R19").Value = "Hortaliza"
